{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Change being applied (per the commit \"solved problem of heretical\n// lowercase letters\" and the accompanying OOXML diff):\n//   1. Capitalise the first \"visualisations\" -> \"Visualisations\" in the\n//      sentence \"At Imperial visualisations, we have several design\n//      principles ...\" (only the occurrence right after \"At Imperial\").\n//   2. Remove the stray \"_GoBack\" bookmark left over near the end of the\n//      document (just after the font-guide hyperlink paragraph).\n\nconst body = context.document.body;\n\n// ---- 1) \"At Imperial visualisations\" -> \"At Imperial Visualisations\" ----\nconst phraseResults = body.search(\"At Imperial visualisations, we have\", {\n  matchCase: true\n});\nphraseResults.load(\"items\");\nawait context.sync();\n\nif (phraseResults.items.length > 0) {\n  // Narrow the search down to the paragraph that contains the phrase so we\n  // only touch the \"v\" immediately following \"At Imperial \" (the word\n  // \"visualisations\"/\"visualisation\" recurs several times later in the\n  // same paragraph and document).\n  const matchedRange = phraseResults.items[0];\n  const paragraphs = matchedRange.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  const vResults = paragraph.search(\"v\", { matchCase: true });\n  vResults.load(\"items\");\n  await context.sync();\n\n  if (vResults.items.length > 0) {\n    // Search results come back in document order, so the first hit is the\n    // \"v\" that starts \"visualisations\" right after \"At Imperial \".\n    vResults.items[0].insertText(\"V\", Word.InsertLocation.replace);\n  }\n}\n\n// ---- 2) Drop the leftover \"_GoBack\" bookmark ----\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Change being applied (per the commit \"solved problem of heretical\n# lowercase letters\" and the accompanying OOXML diff):\n#   1. Capitalise the first \"visualisations\" -> \"Visualisations\" in the\n#      sentence \"At Imperial visualisations, we have several design\n#      principles ...\" (only the occurrence right after \"At Imperial\").\n#   2. Remove the stray \"_GoBack\" bookmark left over near the end of the\n#      document (just after the font-guide hyperlink paragraph).\n\n$d = $word.ActiveDocument\n\n# ---- 1) \"At Imperial visualisations\" -> \"At Imperial Visualisations\" ----\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"At Imperial visualisations, we have\"\n$found = $find.Execute()\nif ($found) {\n    # $rng now spans exactly \"At Imperial visualisations, we have\"; narrow\n    # down to the single \"v\" that starts \"visualisations\" (right after\n    # \"At Imperial \") so only that character is touched, leaving the other\n    # (later, lowercase) occurrences of \"visualisations\" untouched.\n    $vStart = $rng.Start + \"At Imperial \".Length\n    $vRange = $d.Range($vStart, $vStart + 1)\n    if ($vRange.Text -eq \"v\") {\n        $vRange.Text = \"V\"\n    }\n}\n\n# ---- 2) Drop the leftover \"_GoBack\" bookmark ----\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
